$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "PTH24u + dCa"
$ws.Range("A5").Value = "PTH24u + dCorrCa"
$ws.Range("A8").Value = "dPTH + dCa"
$ws.Range("A9").Value = "dPTH + dCorrCa"
